$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.223871366608478
$ws.Range("C2").Value = 0.2914999872252224
$ws.Range("D2").Value = 0.0146276085420709
$ws.Range("F2").Value = 4.190190993022128
$ws.Range("G2").Value = 0.002622544685130013
$ws.Range("I2").Value = 2.547437221437491
$ws.Range("J2").Value = 0.1332338976556851
$ws.Range("L2").Value = 0.3829061384361552
$ws.Range("M2").Value = 0.5025365823670427
$ws.Range("B3").Value = 2.156817395947485
$ws.Range("C3").Value = 0.266344581104903
$ws.Range("D3").Value = 0.01338830295732407
$ws.Range("F3").Value = 4.158147181707932
$ws.Range("G3").Value = 0.002628226200481901
$ws.Range("I3").Value = 2.533548856834258
$ws.Range("J3").Value = 0.1331434208125706
$ws.Range("L3").Value = 0.3818556923943675
$ws.Range("M3").Value = 0.4929297616668791
$ws.Range("B4").Value = 2.117082324608248
$ws.Range("C4").Value = 0.2510732038874721
$ws.Range("D4").Value = 0.01262275032486926
$ws.Range("F4").Value = 4.140288958729741
$ws.Range("G4").Value = 0.002631898310926245
$ws.Range("I4").Value = 2.526026167812574
$ws.Range("J4").Value = 0.1330856449885385
$ws.Range("L4").Value = 0.3813733727305859
$ws.Range("M4").Value = 0.4873288713017701
$ws.Range("B5").Value = 2.101251029088701
$ws.Range("C5").Value = 0.2448933027155249
$ws.Range("D5").Value = 0.01230954806471019
$ws.Range("F5").Value = 4.133466852501755
$ws.Range("G5").Value = 0.002633441061900575
$ws.Range("I5").Value = 2.523212306378369
$ws.Range("J5").Value = 0.1330615369210548
$ws.Range("L5").Value = 0.3812177740218843
$ws.Range("M5").Value = 0.4851213497459028
$ws.Range("B6").Value = 2.098644059776916
$ws.Range("C6").Value = 0.2438697367140605
$ws.Range("D6").Value = 0.01225746444288589
$ws.Range("F6").Value = 4.132361497716474
$ws.Range("G6").Value = 0.002633700037931241
$ws.Range("I6").Value = 2.522760240192
$ws.Range("J6").Value = 0.1330574995876259
$ws.Range("L6").Value = 0.3811944117528512
$ws.Range("M6").Value = 0.4847593161275867
$ws.Range("B7").Value = 2.116867356402793
$ws.Range("C7").Value = 0.2509896848489177
$ws.Range("D7").Value = 0.01261853145954817
$ws.Range("F7").Value = 4.140195112198825
$ws.Range("G7").Value = 0.002631918929295999
$ws.Range("I7").Value = 2.525987201301007
$ws.Range("J7").Value = 0.1330853221482484
$ws.Range("L7").Value = 0.3813711083994065
$ws.Range("M7").Value = 0.4872987966542865
$ws.Range("B8").Value = 2.200452940967637
$ws.Range("C8").Value = 0.2827899599454895
$ws.Range("D8").Value = 0.01420120928597868
$ws.Range("F8").Value = 4.178764233009105
$ws.Range("G8").Value = 0.002624465653888372
$ws.Range("I8").Value = 2.542439365503853
$ws.Range("J8").Value = 0.1332031595001268
$ws.Range("L8").Value = 0.3825102123347719
$ws.Range("M8").Value = 0.4991623602326669
$ws.Range("B9").Value = 2.375779452148493
$ws.Range("C9").Value = 0.3465573184786308
$ws.Range("D9").Value = 0.01727162818077232
$ws.Range("F9").Value = 4.268893361817277
$ws.Range("G9").Value = 0.002611299604335482
$ws.Range("I9").Value = 2.582722657592853
$ws.Range("J9").Value = 0.1334168344353848
$ws.Range("L9").Value = 0.3860330940617445
$ws.Range("M9").Value = 0.5247906172511634
$ws.Range("B10").Value = 2.511596866500156
$ws.Range("C10").Value = 0.3943054336999694
$ws.Range("D10").Value = 0.01951208163445983
$ws.Range("F10").Value = 4.344069920428808
$ws.Range("G10").Value = 0.002602500145324847
$ws.Range("I10").Value = 2.617280159543839
$ws.Range("J10").Value = 0.1335635401227169
$ws.Range("L10").Value = 0.3894061050882272
$ws.Range("M10").Value = 0.5450660496549773
$ws.Range("B11").Value = 2.574916535654211
$ws.Range("C11").Value = 0.416231500345873
$ws.Range("D11").Value = 0.02052904155421942
$ws.Range("F11").Value = 4.380242269709669
$ws.Range("G11").Value = 0.002598684566870488
$ws.Range("I11").Value = 2.634094501042469
$ws.Range("J11").Value = 0.1336281198529337
$ws.Range("L11").Value = 0.391110834741923
$ws.Range("M11").Value = 0.5546053817431371
$ws.Range("B12").Value = 2.599115509775174
$ws.Range("C12").Value = 0.4245644911586055
$ws.Range("D12").Value = 0.02091389969553603
$ws.Range("F12").Value = 4.394225698310152
$ws.Range("G12").Value = 0.002597266479209711
$ws.Range("I12").Value = 2.64062016437714
$ws.Range("J12").Value = 0.1336522700776079
$ws.Range("L12").Value = 0.3917808406490764
$ws.Range("M12").Value = 0.5582631775824964
$ws.Range("B13").Value = 2.593893983989631
$ws.Range("C13").Value = 0.4227684863545846
$ws.Range("D13").Value = 0.02083102313571317
$ws.Range("F13").Value = 4.391201378253356
$ws.Range("G13").Value = 0.002597570700750766
$ws.Range("I13").Value = 2.639207681420999
$ws.Range("J13").Value = 0.1336470823584301
$ws.Range("L13").Value = 0.3916354554222607
$ws.Range("M13").Value = 0.5574733835567187
$ws.Range("B14").Value = 2.57690296460953
$ws.Range("C14").Value = 0.4169164547353716
$ws.Range("D14").Value = 0.02056070847256564
$ws.Range("F14").Value = 4.381386957288157
$ws.Range("G14").Value = 0.002598567363869944
$ws.Range("I14").Value = 2.634628189485227
$ws.Range("J14").Value = 0.1336301127819315
$ws.Range("L14").Value = 0.3911654664529181
$ws.Range("M14").Value = 0.5549053995894724
$ws.Range("B15").Value = 2.566524291603116
$ws.Range("C15").Value = 0.4133358484780274
$ws.Range("D15").Value = 0.02039510369852593
$ws.Range("F15").Value = 4.375412609183485
$ws.Range("G15").Value = 0.002599181333376512
$ws.Range("I15").Value = 2.631843785257331
$ws.Range("J15").Value = 0.1336196789077144
$ws.Range("L15").Value = 0.3908807693890708
$ws.Range("M15").Value = 0.5533383561148995
$ws.Range("B16").Value = 2.507489708647768
$ws.Range("C16").Value = 0.392876690326375
$ws.Range("D16").Value = 0.01944558267963714
$ws.Range("F16").Value = 4.3417458562067
$ws.Range("G16").Value = 0.002602753259132632
$ws.Range("I16").Value = 2.616203412161028
$ws.Range("J16").Value = 0.1335592767422185
$ws.Range("L16").Value = 0.3892981217317981
$ws.Range("M16").Value = 0.5444489914267976
$ws.Range("B17").Value = 2.471667467079953
$ws.Range("C17").Value = 0.3803786109415341
$ws.Range("D17").Value = 0.01886256602719527
$ws.Range("F17").Value = 4.321599341440276
$ws.Range("G17").Value = 0.00260499239513487
$ws.Range("I17").Value = 2.606889537932872
$ws.Range("J17").Value = 0.1335216733045774
$ws.Range("L17").Value = 0.3883708211781851
$ws.Range("M17").Value = 0.5390765896828924
$ws.Range("B18").Value = 2.45120797909658
$ws.Range("C18").Value = 0.3732093275954185
$ws.Range("D18").Value = 0.01852701529541889
$ws.Range("F18").Value = 4.310197303162596
$ws.Range("G18").Value = 0.002606297929471907
$ws.Range("I18").Value = 2.601635344595167
$ws.Range("J18").Value = 0.1334998413999173
$ws.Range("L18").Value = 0.3878534964411671
$ws.Range("M18").Value = 0.5360162616946909
$ws.Range("B19").Value = 2.444305558614019
$ws.Range("C19").Value = 0.3707852262599545
$ws.Range("D19").Value = 0.01841336427322915
$ws.Range("F19").Value = 4.306368607329347
$ws.Range("G19").Value = 0.002606742995281548
$ws.Range("I19").Value = 2.599874004132573
$ws.Range("J19").Value = 0.1334924143911458
$ws.Range("L19").Value = 0.3876810940251261
$ws.Range("M19").Value = 0.5349851936051024
$ws.Range("B20").Value = 2.475465846280031
$ws.Range("C20").Value = 0.3817070532867888
$ws.Range("D20").Value = 0.01892465082513439
$ws.Range("F20").Value = 4.32372473923013
$ws.Range("G20").Value = 0.002604752210655092
$ws.Range("I20").Value = 2.607870357778253
$ws.Range("J20").Value = 0.1335256972659833
$ws.Range("L20").Value = 0.3884678746803019
$ws.Range("M20").Value = 0.5396454133489996
$ws.Range("B21").Value = 2.581887632728069
$ws.Range("C21").Value = 0.4186345178273427
$ws.Range("D21").Value = 0.02064011239435359
$ws.Range("F21").Value = 4.384261920627267
$ws.Range("G21").Value = 0.002598273894012425
$ws.Range("I21").Value = 2.635968988302352
$ws.Range("J21").Value = 0.1336351053836857
$ws.Range("L21").Value = 0.3913028499598568
$ws.Range("M21").Value = 0.5556584447716162
$ws.Range("B22").Value = 2.652729883675818
$ws.Range("C22").Value = 0.4429442599279128
$ws.Range("D22").Value = 0.02175988906893878
$ws.Range("F22").Value = 4.425492889644204
$ws.Range("G22").Value = 0.002594196017507544
$ws.Range("I22").Value = 2.655257091079747
$ws.Range("J22").Value = 0.1337048363998248
$ws.Range("L22").Value = 0.3932982391574598
$ws.Range("M22").Value = 0.5663888666301276
$ws.Range("B23").Value = 2.614801923750065
$ws.Range("C23").Value = 0.4299534472287974
$ws.Range("D23").Value = 0.02116234300650888
$ws.Range("F23").Value = 4.403334046111297
$ws.Range("G23").Value = 0.002596358225166302
$ws.Range("I23").Value = 2.644877744284869
$ws.Range("J23").Value = 0.1336677800778521
$ws.Range("L23").Value = 0.3922202271418342
$ws.Range("M23").Value = 0.56063758225973
$ws.Range("B24").Value = 2.473748177965376
$ws.Range("C24").Value = 0.3811064146176477
$ws.Range("D24").Value = 0.01889658343604594
$ws.Range("F24").Value = 4.322763284932762
$ws.Range("G24").Value = 0.00260486074141912
$ws.Range("I24").Value = 2.607426616247622
$ws.Range("J24").Value = 0.1335238786970843
$ws.Range("L24").Value = 0.3884239475987528
$ws.Range("M24").Value = 0.5393881599069275
$ws.Range("B25").Value = 2.3271222892065
$ws.Range("C25").Value = 0.3291517766342338
$ws.Range("D25").Value = 0.01644405284012862
$ws.Range("F25").Value = 4.242946692876416
$ws.Range("G25").Value = 0.002614707215404617
$ws.Range("I25").Value = 2.570958909192115
$ws.Range("J25").Value = 0.1333608591908488
$ws.Range("L25").Value = 0.384942177719708
$ws.Range("M25").Value = 0.5176039066524254
